$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.18070387840271
$ws.Range("B1").Value = 2.394456148147583
$ws.Range("C1").Value = 3.651184320449829
$ws.Range("D1").Value = 2.057796478271484
$ws.Range("E1").Value = 1.204771399497986
